# Property disclosure workbook update (issue #5): property boat&car done.
# The 汽車 (Car) sheet gets a proper header row (field-name labels instead of
# duplicated row-2 data), a new "capacity" column, and the standard trailing
# metadata columns (property_category / category / date / legislator_name /
# legislator_id / source_file / index) that the other sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Header row (row 1): replace the old mirrored-data header with the
# canonical field-name labels, and extend it with the metadata columns. ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data rows (2 = BMW, 3 = Audi A5): columns A-G already hold the right
# values (name/capacity/owner/register_date/register_reason/acquire_value),
# so only the new trailing metadata columns need to be filled in. ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2011-11-22"
$ws.Range("K2").Value = "陳根德"
$ws.Range("L2").Value = 833
$ws.Range("M2").Value = "tmpa3b61"
$ws.Range("N2").Value = 31

$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2011-11-22"
$ws.Range("K3").Value = "陳根德"
$ws.Range("L3").Value = 833
$ws.Range("M3").Value = "tmpa3b61"
$ws.Range("N3").Value = 32

# --- Formatting: copy the existing header / data-row look onto the new
# columns so H:N match the B:G styling (bold+border header, bordered data). ---
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("H3:N3").PasteSpecial(-4122)
